$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "Contact Name" column header to "Resident Name".
$ws.Range("F1").Value = "Resident Name"

# Update the cell comments to reflect the new "resident" terminology.
$null = $ws.Range("F1").Comment.Text("This is a mandatory field")
$null = $ws.Range("G1").Comment.Text("Email of the resident. All email around tenant access and secret changes will be sent to this address.")
$null = $ws.Range("H1").Comment.Text("Phone number of the resident. Format +<country code><phone number>")

# Move the active selection, matching the author's last cursor position.
$null = $ws.Range("G22").Select()
